$d = $word.ActiveDocument

$d.Content.Find.Execute("772÷7=110, 2", $true, $false, $false, $false, $false, $true, 1, $false, "882÷5=176, 2", 2) | Out-Null
$d.Content.Find.Execute("775÷2=387, 1", $true, $false, $false, $false, $false, $true, 1, $false, "194÷8=24, 2", 2) | Out-Null
$d.Content.Find.Execute("915÷4=228, 3", $true, $false, $false, $false, $false, $true, 1, $false, "882÷3=294, 0", 2) | Out-Null
$d.Content.Find.Execute("737÷9=81, 8", $true, $false, $false, $false, $false, $true, 1, $false, "434÷9=48, 2", 2) | Out-Null
$d.Content.Find.Execute("761÷4=190, 1", $true, $false, $false, $false, $false, $true, 1, $false, "524÷2=262, 0", 2) | Out-Null
$d.Content.Find.Execute("963÷2=481, 1", $true, $false, $false, $false, $false, $true, 1, $false, "960÷5=192, 0", 2) | Out-Null
$d.Content.Find.Execute("431÷8=53, 7", $true, $false, $false, $false, $false, $true, 1, $false, "766÷3=255, 1", 2) | Out-Null
$d.Content.Find.Execute("183÷9=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "871÷8=108, 7", 2) | Out-Null
$d.Content.Find.Execute("503÷8=62, 7", $true, $false, $false, $false, $false, $true, 1, $false, "728÷8=91, 0", 2) | Out-Null
$d.Content.Find.Execute("623÷5=124, 3", $true, $false, $false, $false, $false, $true, 1, $false, "250÷4=62, 2", 2) | Out-Null
$d.Content.Find.Execute("925÷8=115, 5", $true, $false, $false, $false, $false, $true, 1, $false, "424÷6=70, 4", 2) | Out-Null
$d.Content.Find.Execute("680÷5=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "363÷9=40, 3", 2) | Out-Null
$d.Content.Find.Execute("131÷9=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "937÷8=117, 1", 2) | Out-Null
$d.Content.Find.Execute("934÷5=186, 4", $true, $false, $false, $false, $false, $true, 1, $false, "260÷7=37, 1", 2) | Out-Null
$d.Content.Find.Execute("499÷5=99, 4", $true, $false, $false, $false, $false, $true, 1, $false, "723÷6=120, 3", 2) | Out-Null
$d.Content.Find.Execute("348÷6=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "871÷4=217, 3", 2) | Out-Null
$d.Content.Find.Execute("717÷3=239, 0", $true, $false, $false, $false, $false, $true, 1, $false, "564÷6=94, 0", 2) | Out-Null
$d.Content.Find.Execute("482÷3=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "269÷5=53, 4", 2) | Out-Null
$d.Content.Find.Execute("584÷5=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "753÷7=107, 4", 2) | Out-Null
$d.Content.Find.Execute("141÷2=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "133÷4=33, 1", 2) | Out-Null
$d.Content.Find.Execute("599÷2=299, 1", $true, $false, $false, $false, $false, $true, 1, $false, "699÷5=139, 4", 2) | Out-Null
$d.Content.Find.Execute("389÷3=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "915÷6=152, 3", 2) | Out-Null
$d.Content.Find.Execute("987÷9=109, 6", $true, $false, $false, $false, $false, $true, 1, $false, "339÷8=42, 3", 2) | Out-Null
$d.Content.Find.Execute("729÷2=364, 1", $true, $false, $false, $false, $false, $true, 1, $false, "203÷2=101, 1", 2) | Out-Null
$d.Content.Find.Execute("187÷8=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "978÷5=195, 3", 2) | Out-Null
